# edit.ps1 -- applies the Record.docx diff via Word COM-interop (iron_native)
#
# Summary of changes:
#  1. Paragraph "call_gamma*current_call_position-put_gamma" gets spell-check
#     proofErr markup added around its tokens (text unchanged).
#  2. Paragraph "money_account + ... basemodel.callPrice(...) ..." gets
#     proofErr markup added (text unchanged).
#  3. Paragraph "current_stock_position + current_call_position * call_delta -
#     current_delta" gets proofErr markup added (text unchanged).
#  4. Two brand-new paragraphs are inserted right after that paragraph (with
#     one blank separator paragraph before them):
#       "stock_position-put_delta+call_position*call_price"
#       "-base.Gamma(St[k], 0.25, day)+base.Gamma(St[k], 0.5, day)*call_position"
#     the second one carries a <w:rFonts w:hint="eastAsia"/> run-property on
#     its paragraph mark.
#  5. The "cVar " run is split into a proofErr-wrapped "cVar" run plus a
#     separate " " run (text unchanged).
#  6. The "DG Move based cvar = " run is split into a plain "DG Move based "
#     run plus a proofErr-wrapped "cvar" run plus a plain " = " run (text
#     unchanged).

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Replace-ParagraphXml($paragraph, [string]$innerXml) {
    # Replace a whole paragraph (incl. its paragraph mark) with fresh XML,
    # in place -- does not change the paragraph count.
    $full = $paragraph.Range
    $xml = '<w:p ' + $wNs + '>' + $innerXml + '</w:p>'
    $full.InsertXML($xml)
}

# ---------------------------------------------------------------------
# 1) call_gamma*current_call_position-put_gamma
# ---------------------------------------------------------------------
$p = $d.Paragraphs.Item(4)
$inner = '<w:proofErr w:type="spellStart"/>' +
         '<w:r><w:t>call_gamma</w:t></w:r>' +
         '<w:proofErr w:type="spellEnd"/>' +
         '<w:r><w:t>*</w:t></w:r>' +
         '<w:proofErr w:type="spellStart"/>' +
         '<w:r><w:t>current_call_position-put_gamma</w:t></w:r>' +
         '<w:proofErr w:type="spellEnd"/>'
Replace-ParagraphXml $p $inner

# ---------------------------------------------------------------------
# 2) money_account + current_call_position*basemodel.callPrice(...)...
# ---------------------------------------------------------------------
$p = $d.Paragraphs.Item(6)
$inner = '<w:proofErr w:type="spellStart"/>' +
         '<w:r><w:t>money_account</w:t></w:r>' +
         '<w:proofErr w:type="spellEnd"/>' +
         '<w:r><w:t xml:space="preserve"> + current_call_position*</w:t></w:r>' +
         '<w:proofErr w:type="gramStart"/>' +
         '<w:r><w:t>basemodel.callPrice</w:t></w:r>' +
         '<w:proofErr w:type="gramEnd"/>' +
         '<w:r><w:t>(St[k+1],0.5,0.25/91*(k+1))+current_stock_position*St[k+1]-basemodel.putPrice(St[k+1], (k+1)*0.25/91)</w:t></w:r>'
Replace-ParagraphXml $p $inner

# ---------------------------------------------------------------------
# 3) current_stock_position + current_call_position * call_delta - current_delta
# ---------------------------------------------------------------------
$p = $d.Paragraphs.Item(11)
$inner = '<w:proofErr w:type="spellStart"/>' +
         '<w:r><w:t>current_stock_position</w:t></w:r>' +
         '<w:proofErr w:type="spellEnd"/>' +
         '<w:r><w:t xml:space="preserve"> + </w:t></w:r>' +
         '<w:proofErr w:type="spellStart"/>' +
         '<w:r><w:t>current_call_position</w:t></w:r>' +
         '<w:proofErr w:type="spellEnd"/>' +
         '<w:r><w:t xml:space="preserve"> * </w:t></w:r>' +
         '<w:proofErr w:type="spellStart"/>' +
         '<w:r><w:t>call_delta</w:t></w:r>' +
         '<w:proofErr w:type="spellEnd"/>' +
         '<w:r><w:t xml:space="preserve"> - </w:t></w:r>' +
         '<w:proofErr w:type="spellStart"/>' +
         '<w:r><w:t>current_delta</w:t></w:r>' +
         '<w:proofErr w:type="spellEnd"/>'
Replace-ParagraphXml $p $inner

# ---------------------------------------------------------------------
# 4) insert 3 new paragraphs right after it: blank, then the two new lines
# ---------------------------------------------------------------------
$anchor = $d.Paragraphs.Item(11).Range
$anchor.InsertParagraphAfter()
$anchor.InsertParagraphAfter()
$anchor.InsertParagraphAfter()

# paragraph 12 is the blank separator -- clear the stray empty run Word
# leaves behind so it serialises as a truly empty <w:p/>
$pA = $d.Paragraphs.Item(12)
Replace-ParagraphXml $pA ''

# fill 13 and 14
$pB = $d.Paragraphs.Item(13)
$innerB = '<w:proofErr w:type="spellStart"/>' +
          '<w:r><w:t>stock_position-put_delta+call_position</w:t></w:r>' +
          '<w:proofErr w:type="spellEnd"/>' +
          '<w:r><w:t>*</w:t></w:r>' +
          '<w:proofErr w:type="spellStart"/>' +
          '<w:r><w:t>call_price</w:t></w:r>' +
          '<w:proofErr w:type="spellEnd"/>'
Replace-ParagraphXml $pB $innerB

$pC = $d.Paragraphs.Item(14)
$innerC = '<w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr>' +
          '<w:r><w:t>-</w:t></w:r>' +
          '<w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/>' +
          '<w:r><w:t>base.Gamma</w:t></w:r>' +
          '<w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/>' +
          '<w:r><w:t>(St[k], 0.25, day)+</w:t></w:r>' +
          '<w:proofErr w:type="spellStart"/>' +
          '<w:r><w:t>base.Gamma</w:t></w:r>' +
          '<w:proofErr w:type="spellEnd"/>' +
          '<w:r><w:t>(St[k], 0.5, day)*</w:t></w:r>' +
          '<w:proofErr w:type="spellStart"/>' +
          '<w:r><w:t>call_position</w:t></w:r>' +
          '<w:proofErr w:type="spellEnd"/>'
Replace-ParagraphXml $pC $innerC

# ---------------------------------------------------------------------
# 5) cVar -0.3235655547073288  ->  split "cVar " run into "cVar" + " "
#    (the remaining three paragraphs below shifted down by +3)
# ---------------------------------------------------------------------
$p = $d.Paragraphs.Item(48)
$inner = '<w:pPr><w:ind w:firstLineChars="100" w:firstLine="210"/></w:pPr>' +
         '<w:proofErr w:type="spellStart"/>' +
         '<w:r><w:t>cVar</w:t></w:r>' +
         '<w:proofErr w:type="spellEnd"/>' +
         '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
         '<w:r><w:t>-0.3235655547073288</w:t></w:r>'
Replace-ParagraphXml $p $inner

# ---------------------------------------------------------------------
# 6) DG Move based cvar = -0.26538250626172377
#    -> split into "DG Move based " + proofErr("cvar") + " = "
# ---------------------------------------------------------------------
$p = $d.Paragraphs.Item(55)
$inner = '<w:r><w:t xml:space="preserve">DG Move based </w:t></w:r>' +
         '<w:proofErr w:type="spellStart"/>' +
         '<w:r><w:t>cvar</w:t></w:r>' +
         '<w:proofErr w:type="spellEnd"/>' +
         '<w:r><w:t xml:space="preserve"> = </w:t></w:r>' +
         '<w:r><w:t>-0.26538250626172377</w:t></w:r>'
Replace-ParagraphXml $p $inner

Write-Output "done"
